# Extend the "working_hours" sheet with two more time-tracking entries
# (2014-02-18 and 2014-02-20), push the blank separator row and the
# summary rows down accordingly, and re-create the time-spent formula
# in column F as a single shared formula spanning F2:F9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert two new rows at 8:9, pushing the existing blank
# separator row (old row 8) and the three summary rows (old rows 9-11)
# down to rows 10-13. Excel copies row 8's formatting into the new rows
# and auto-adjusts every formula reference (SUM/F9/F10 etc.) that spans
# the insertion point.
$ws.Rows("8:9").Insert()

# New data row for 2014-02-18, 20:00-21:00
$ws.Range("A8").Value = 2014
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 18
$ws.Range("D8").Value = 0.83333333333333337
$ws.Range("E8").Value = 0.875

# New data row for 2014-02-20, 11:00-12:00
$ws.Range("A9").Value = 2014
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 20
$ws.Range("D9").Value = 0.45833333333333331
$ws.Range("E9").Value = 0.5

# Re-write the "time spent [min]" formula across the full data range as
# one shared formula (F2:F9) instead of eight independent copies.
$ws.Range("F2:F9").Formula = "=(E2-D2)*24*60"

# Match the author's final cursor position.
[void]$ws.Range("D10").Select()
